# Update stock quantities in the BOM "Remarks_customer" column (column M)
# to reflect the latest stock counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "M2"  = "Stock: 14520"
    "M3"  = "Stock: 66390"
    "M4"  = "Stock: 6860"
    "M6"  = "Stock: 132555"
    "M7"  = "Stock: 46486"
    "M13" = "Stock: 341980"
    "M14" = "Stock: 9322"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
